# Auto-relevant edit script: updates GUID/hash/timestamp strings produced by a fresh
# handback-status report generation run.
$wb = $excel.ActiveWorkbook

# ---- Sheet: Overview ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/23052a6ce28d398fb9ac145ea5f60809c36aa1bc/e2e/26c045ad-ebac-414f-90e5-58c5441bf91d.md", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.md") | Out-Null
$ws.Range("A3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/23052a6ce28d398fb9ac145ea5f60809c36aa1bc/e2e/8d7f325e-639e-470a-8973-9a04833e45c6.md", "", "", "ffffdca2239e-f0b5-4af2-8abf-4cec5d47727c.md") | Out-Null

# ---- Sheet: zh-cn ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/23052a6ce28d398fb9ac145ea5f60809c36aa1bc/e2e/26c045ad-ebac-414f-90e5-58c5441bf91d.md", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.md") | Out-Null
$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/23052a6ce28d398fb9ac145ea5f60809c36aa1bc/e2e/26c045ad-ebac-414f-90e5-58c5441bf91d.md", "", "", ".md") | Out-Null
$ws.Range("D2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/74997808139033eb07e559b7ca83db338b08ab55/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/26c045ad-ebac-414f-90e5-58c5441bf91d.5fb125723dc015112c84191c8784859b82433665.zh-cn.xlf", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.c0e5cbbd5c3502ce52a22c5403ff31242e293cbc.zh-cn.xlf") | Out-Null
$ws.Range("F2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/794552b5aa6218a482673f57515ece4253f2d54e/e2e/26c045ad-ebac-414f-90e5-58c5441bf91d.md", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.md") | Out-Null
$ws.Range("G2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d91e1bfeffbf9a219983f1d4b19922acf18fc192/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/26c045ad-ebac-414f-90e5-58c5441bf91d.5fb125723dc015112c84191c8784859b82433665.zh-cn.xlf", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.c0e5cbbd5c3502ce52a22c5403ff31242e293cbc.zh-cn.xlf") | Out-Null
$ws.Range("A3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/23052a6ce28d398fb9ac145ea5f60809c36aa1bc/e2e/8d7f325e-639e-470a-8973-9a04833e45c6.md", "", "", "ffffdca2239e-f0b5-4af2-8abf-4cec5d47727c.md") | Out-Null
$ws.Range("B3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/23052a6ce28d398fb9ac145ea5f60809c36aa1bc/e2e/8d7f325e-639e-470a-8973-9a04833e45c6.md", "", "", ".md") | Out-Null
$ws.Range("D3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/74997808139033eb07e559b7ca83db338b08ab55/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8d7f325e-639e-470a-8973-9a04833e45c6.991a236844311b5d2260d6272d2d6e196e409901.zh-cn.xlf", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.c0e5cbbd5c3502ce52a22c5403ff31242e293cbc.zh-cn.xlf") | Out-Null
$ws.Range("F3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/794552b5aa6218a482673f57515ece4253f2d54e/e2e/8d7f325e-639e-470a-8973-9a04833e45c6.md", "", "", "ffffdca2239e-f0b5-4af2-8abf-4cec5d47727c.md") | Out-Null
$ws.Range("G3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d91e1bfeffbf9a219983f1d4b19922acf18fc192/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8d7f325e-639e-470a-8973-9a04833e45c6.991a236844311b5d2260d6272d2d6e196e409901.zh-cn.xlf", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.c0e5cbbd5c3502ce52a22c5403ff31242e293cbc.zh-cn.xlf") | Out-Null

# ---- Sheet: de-de ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/23052a6ce28d398fb9ac145ea5f60809c36aa1bc/e2e/26c045ad-ebac-414f-90e5-58c5441bf91d.md", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.md") | Out-Null
$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/23052a6ce28d398fb9ac145ea5f60809c36aa1bc/e2e/26c045ad-ebac-414f-90e5-58c5441bf91d.md", "", "", ".md") | Out-Null
$ws.Range("D2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/78867b9e11d27bd17a8a0ca800acb4febbf1be59/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/26c045ad-ebac-414f-90e5-58c5441bf91d.5fb125723dc015112c84191c8784859b82433665.de-de.xlf", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.c0e5cbbd5c3502ce52a22c5403ff31242e293cbc.de-de.xlf") | Out-Null
$ws.Range("F2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fd9da5bef4670c1046c848ee625c6ce81cf3ad9b/e2e/26c045ad-ebac-414f-90e5-58c5441bf91d.md", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.md") | Out-Null
$ws.Range("G2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1e280a8c4aa7482c6271ec0151a4995d8a3c65d4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/26c045ad-ebac-414f-90e5-58c5441bf91d.5fb125723dc015112c84191c8784859b82433665.de-de.xlf", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.c0e5cbbd5c3502ce52a22c5403ff31242e293cbc.de-de.xlf") | Out-Null
$ws.Range("A3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/23052a6ce28d398fb9ac145ea5f60809c36aa1bc/e2e/8d7f325e-639e-470a-8973-9a04833e45c6.md", "", "", "ffffdca2239e-f0b5-4af2-8abf-4cec5d47727c.md") | Out-Null
$ws.Range("B3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/23052a6ce28d398fb9ac145ea5f60809c36aa1bc/e2e/8d7f325e-639e-470a-8973-9a04833e45c6.md", "", "", ".md") | Out-Null
$ws.Range("D3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/78867b9e11d27bd17a8a0ca800acb4febbf1be59/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8d7f325e-639e-470a-8973-9a04833e45c6.991a236844311b5d2260d6272d2d6e196e409901.de-de.xlf", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.c0e5cbbd5c3502ce52a22c5403ff31242e293cbc.de-de.xlf") | Out-Null
$ws.Range("F3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fd9da5bef4670c1046c848ee625c6ce81cf3ad9b/e2e/8d7f325e-639e-470a-8973-9a04833e45c6.md", "", "", "ffffdca2239e-f0b5-4af2-8abf-4cec5d47727c.md") | Out-Null
$ws.Range("G3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1e280a8c4aa7482c6271ec0151a4995d8a3c65d4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8d7f325e-639e-470a-8973-9a04833e45c6.991a236844311b5d2260d6272d2d6e196e409901.de-de.xlf", "", "", "597a00e2-acda-459d-b80a-0c78127569f3.c0e5cbbd5c3502ce52a22c5403ff31242e293cbc.de-de.xlf") | Out-Null

# ---- Date/time text cells (no hyperlink, just shared-string text) ----
$wsZh = $wb.Worksheets.Item(2)
$wsZh.Range("E2").Value = "2016-03-20 18:50:01"
$wsZh.Range("H2").Value = "2016-03-20 18:50:21"
$wsZh.Range("E3").Value = "2016-03-20 18:50:01"
$wsZh.Range("H3").Value = "2016-03-20 18:50:21"

$wsDe = $wb.Worksheets.Item(3)
$wsDe.Range("E2").Value = "2016-03-20 18:50:06"
$wsDe.Range("H2").Value = "2016-03-20 18:50:27"
$wsDe.Range("E3").Value = "2016-03-20 18:50:06"
$wsDe.Range("H3").Value = "2016-03-20 18:50:27"
